$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1000
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = $null

$ws.Range("H20").Value = 33407
$ws.Range("I20").Value = 33407
$ws.Range("K20").Value = 33407
$ws.Range("M20").Value = -33177

$ws.Range("H21").Value = 42072.93
$ws.Range("I21").Value = 43000
$ws.Range("J21").Value = 41820.09
$ws.Range("K21").Value = 43000
$ws.Range("L21").Value = 41820.09
$ws.Range("M21").Value = -42532
$ws.Range("N21").Value = -42756.09

$ws.Range("H23").Value = 42072.93
$ws.Range("I23").Value = 43000
$ws.Range("J23").Value = 41820.09
$ws.Range("K23").Value = 43000
$ws.Range("L23").Value = 41820.09
$ws.Range("M23").Value = -42766
$ws.Range("N23").Value = -42288.09

$ws.Range("H34").Value = 17843.111
$ws.Range("I34").Value = 17843.111
$ws.Range("K34").Value = 17843.111
$ws.Range("M34").Value = -17640.111

$ws.Range("H35").Value = 33407
$ws.Range("I35").Value = 33407
$ws.Range("K35").Value = 33407
$ws.Range("M35").Value = -33028

$ws.Range("H36").Value = 17843.111
$ws.Range("I36").Value = 17843.111
$ws.Range("K36").Value = 17843.111
$ws.Range("M36").Value = -17128.111

$ws.Range("H45").Value = 6523.3335
$ws.Range("I45").Value = 17093.334
$ws.Range("J45").Value = 3000
$ws.Range("K45").Value = 51280.00199999999
$ws.Range("L45").Value = 9000
$ws.Range("M45").Value = -51088.00199999999
$ws.Range("N45").Value = -9384

$ws.Range("H116").Value = 25001662
$ws.Range("I116").Value = 28573058
$ws.Range("J116").Value = 1900
$ws.Range("K116").Value = 28573058
$ws.Range("L116").Value = 1900
$ws.Range("M116").Value = -28569616
$ws.Range("N116").Value = -8784

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 3458.1667
$ws.Range("I5").Value = 3458.1667
$ws.Range("K5").Value = 3458.1667
$ws.Range("M5").Value = -3346.1667

$ws.Range("H21").Value = 17438.334
$ws.Range("I21").Value = 17438.334
$ws.Range("K21").Value = 17438.334
$ws.Range("M21").Value = -17064.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 3458.1667
$ws.Range("I4").Value = 3458.1667
$ws.Range("K4").Value = 3458.1667
$ws.Range("M4").Value = -3343.1667

$ws.Range("H76").Value = 25833.334
$ws.Range("J76").Value = 25833.334
$ws.Range("L76").Value = 25833.334
$ws.Range("N76").Value = -26463.334

$ws.Range("H79").Value = 25833.334
$ws.Range("J79").Value = 25833.334
$ws.Range("L79").Value = 25833.334
$ws.Range("N79").Value = -28017.334

$ws.Range("H92").Value = 183130.33
$ws.Range("J92").Value = 183130.33
$ws.Range("L92").Value = 183130.33
$ws.Range("N92").Value = -188122.33

$ws.Range("H134").Value = 3513.5
$ws.Range("I134").Value = 3782.4
$ws.Range("J134").Value = 3321.4285
$ws.Range("K134").Value = 11347.2
$ws.Range("L134").Value = 9964.2855
$ws.Range("M134").Value = -8812.2
$ws.Range("N134").Value = -15034.2855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 23978.715
$ws.Range("I41").Value = 1100.5
$ws.Range("J41").Value = 33130
$ws.Range("K41").Value = 1100.5
$ws.Range("L41").Value = 33130
$ws.Range("M41").Value = -672.5
$ws.Range("N41").Value = -33986

$ws.Range("H59").Value = 39958.332
$ws.Range("J59").Value = 44937.5
$ws.Range("L59").Value = 44937.5
$ws.Range("N59").Value = -47227.5

$ws.Range("H60").Value = 17114.285
$ws.Range("I60").Value = 11666.667
$ws.Range("J60").Value = 21200
$ws.Range("K60").Value = 11666.667
$ws.Range("L60").Value = 21200
$ws.Range("M60").Value = -11155.667
$ws.Range("N60").Value = -22222

$ws.Range("H74").Value = 34000
$ws.Range("J74").Value = 34000
$ws.Range("L74").Value = 34000
$ws.Range("N74").Value = -35748

$ws.Range("H77").Value = 34000
$ws.Range("J77").Value = 34000
$ws.Range("L77").Value = 102000
$ws.Range("N77").Value = -110736

$ws.Range("H88").Value = 17177.428
$ws.Range("J88").Value = 17177.428
$ws.Range("L88").Value = 17177.428
$ws.Range("N88").Value = -17989.428

$ws.Range("H91").Value = 17177.428
$ws.Range("J91").Value = 17177.428
$ws.Range("L91").Value = 17177.428
$ws.Range("N91").Value = -19985.428

$ws.Range("H92").Value = 32601
$ws.Range("J92").Value = 32601
$ws.Range("L92").Value = 32601
$ws.Range("N92").Value = -37593

$ws.Range("H96").Value = 22173.5
$ws.Range("J96").Value = 22173.5
$ws.Range("L96").Value = 22173.5
$ws.Range("N96").Value = -27665.5

$ws.Range("H106").Value = 39800
$ws.Range("J106").Value = 39800
$ws.Range("L106").Value = 39800
$ws.Range("N106").Value = -42324

$ws.Range("H122").Value = 6103
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = $null

$ws.Range("H123").Value = 57780
$ws.Range("J123").Value = 57780
$ws.Range("L123").Value = 57780
$ws.Range("N123").Value = -67580

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 54111.11
$ws.Range("J37").Value = 54111.11
$ws.Range("L37").Value = 162333.33
$ws.Range("N37").Value = -162557.33

$ws.Range("H68").Value = 223029.06
$ws.Range("I68").Value = 303601.06
$ws.Range("J68").Value = 1456.0834
$ws.Range("K68").Value = 910803.1799999999
$ws.Range("L68").Value = 4368.2502
$ws.Range("M68").Value = -909992.1799999999
$ws.Range("N68").Value = -5990.2502

$ws.Range("H71").Value = 223029.06
$ws.Range("I71").Value = 303601.06
$ws.Range("J71").Value = 1456.0834
$ws.Range("K71").Value = 2732409.54
$ws.Range("L71").Value = 13104.7506
$ws.Range("M71").Value = -2728353.54
$ws.Range("N71").Value = -21216.7506

$ws.Range("H112").Value = 5208.8237
$ws.Range("J112").Value = 5821.4287
$ws.Range("L112").Value = 17464.2861
$ws.Range("N112").Value = -19680.2861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 24225.445
$ws.Range("J5").Value = 24753.625
$ws.Range("L5").Value = 24753.625
$ws.Range("N5").Value = -24977.625

$ws.Range("H35").Value = 18500
$ws.Range("I35").Value = 18500
$ws.Range("K35").Value = 18500
$ws.Range("M35").Value = -18202

$ws.Range("H122").Value = 2528.889
$ws.Range("I122").Value = 2576.6667
$ws.Range("J122").Value = 2433.3333
$ws.Range("K122").Value = 7730.000100000001
$ws.Range("L122").Value = 7299.999899999999
$ws.Range("M122").Value = -5280.000100000001
$ws.Range("N122").Value = -12199.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1006.93335
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 1009.4545
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1009.4545
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -1599.4545

$ws.Range("H27").Value = 1006.93335
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 1009.4545
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 1009.4545
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -1223.4545

$ws.Range("H56").Value = 55000
$ws.Range("I56").Value = 50000
$ws.Range("J56").Value = 60000
$ws.Range("K56").Value = 50000
$ws.Range("L56").Value = 60000
$ws.Range("M56").Value = -49309
$ws.Range("N56").Value = -61382

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 3500
$ws.Range("J25").Value = 3500
$ws.Range("L25").Value = 3500
$ws.Range("N25").Value = -4086

$ws.Range("H26").Value = 62863
$ws.Range("I26").Value = 50000
$ws.Range("J26").Value = 68008.2
$ws.Range("K26").Value = 50000
$ws.Range("L26").Value = 68008.2
$ws.Range("M26").Value = -49707
$ws.Range("N26").Value = -68594.2

$ws.Range("H37").Value = 59514.5
$ws.Range("I37").Value = 50000
$ws.Range("K37").Value = 50000
$ws.Range("M37").Value = -49797
